$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells hold numeric-looking values stored as TEXT (shared
# strings), e.g. "7.8". Plain `.Value = "7.82"` would let Excel's
# type-inference re-cast the string as a Number, which would change the
# cell's stored type. Forcing the NumberFormat to Text ("@") before the
# assignment keeps the new value as text too, and resetting the cell
# Style back to "Normal" afterwards avoids leaving a stray text-format
# style applied to the cell (matching the original formatting).
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 13: Enterprises density (per 1000 people)
Set-TextValue "B13" "7.82"
Set-TextValue "C13" "0.77"
Set-TextValue "D13" "8.59"

# Row 14: Employment (% of total)  -- C14 ("33") is unchanged
Set-TextValue "B14" "28.02"
Set-TextValue "D14" "61.02"

# Row 16: Enterprises (% of total)
Set-TextValue "B16" "90.61"
Set-TextValue "C16" "8.96"
Set-TextValue "D16" "99.57"
